# Fixes for VRelay and Upgrade Verification Script
# Update ResultProd (col A) from Fail -> Pass and DateProd (col B) to the
# latest verification run timestamp for the rows that were re-verified.

$wb = $excel.ActiveWorkbook

# --- Sheet: CC-Payments-Sale ---
$ws1 = $wb.Worksheets.Item("CC-Payments-Sale")
$ws1.Range("A2").Value = "Pass"
$ws1.Range("B2").Value = "Fri Aug 22 23:00:26 IST 2025"

# --- Sheet: CC-Payments-Auth ---
$ws2 = $wb.Worksheets.Item("CC-Payments-Auth")
$ws2.Range("A2").Value = "Pass"
$ws2.Range("B2").Value = "Fri Aug 22 22:48:17 IST 2025"
$ws2.Range("A3").Value = "Pass"
$ws2.Range("B3").Value = "Fri Aug 22 22:49:03 IST 2025"
$ws2.Range("A4").Value = "Pass"
$ws2.Range("B4").Value = "Fri Aug 22 22:50:01 IST 2025"
$ws2.Range("A5").Value = "Pass"
$ws2.Range("B5").Value = "Fri Aug 22 22:50:55 IST 2025"
$ws2.Range("A6").Value = "Pass"
$ws2.Range("B6").Value = "Fri Aug 22 22:51:42 IST 2025"
$ws2.Range("A7").Value = "Pass"
$ws2.Range("B7").Value = "Fri Aug 22 22:52:29 IST 2025"

# --- Sheet: ACH-Payments-Debit ---
$ws3 = $wb.Worksheets.Item("ACH-Payments-Debit")
$ws3.Range("A2").Value = "Pass"
$ws3.Range("B2").Value = "Fri Aug 22 22:53:20 IST 2025"
$ws3.Range("A3").Value = "Pass"
$ws3.Range("B3").Value = "Fri Aug 22 22:54:06 IST 2025"
$ws3.Range("A4").Value = "Pass"
$ws3.Range("B4").Value = "Fri Aug 22 22:54:56 IST 2025"
$ws3.Range("A5").Value = "Pass"
$ws3.Range("B5").Value = "Fri Aug 22 22:55:43 IST 2025"
$ws3.Range("A6").Value = "Pass"
$ws3.Range("B6").Value = "Fri Aug 22 22:56:26 IST 2025"
$ws3.Range("A7").Value = "Pass"
$ws3.Range("B7").Value = "Fri Aug 22 22:57:13 IST 2025"
$ws3.Range("A8").Value = "Pass"
$ws3.Range("B8").Value = "Fri Aug 22 22:58:03 IST 2025"
$ws3.Range("A9").Value = "Pass"
$ws3.Range("B9").Value = "Fri Aug 22 22:58:52 IST 2025"
$ws3.Range("A10").Value = "Pass"
$ws3.Range("B10").Value = "Fri Aug 22 22:59:33 IST 2025"
